# Book1.xlsx logic update (kho/lo hang pricing correction):
#  - C2 and C3 (price column) corrected from 1000 -> 13000
#  - Active selection on Sheet1 moved from H6 to P6

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 13000
$ws.Range("C3").Value = 13000

$ws.Range("P6").Select()
